# Update column G ("K") values for rows 2-20 with regenerated s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 2
    9  = 3
    10 = 0
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
